$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the DATA_HORA_ANALISE timestamp string in A2
$ws.Range("A2").Value = "2025-05-28 10:00:08"

# Update the metric values in row 2 (B2:AW2)
$ws.Range("B2").Value = 2540
$ws.Range("C2").Value = 1432
$ws.Range("D2").Value = 56.37795275590551
$ws.Range("E2").Value = 342
$ws.Range("F2").Value = 13.46456692913386
$ws.Range("G2").Value = 1719
$ws.Range("H2").Value = 67.6771653543307
$ws.Range("I2").Value = 804
$ws.Range("J2").Value = 31.65354330708661
$ws.Range("K2").Value = 339953.95
$ws.Range("L2").Value = 17
$ws.Range("M2").Value = 0.6692913385826772
$ws.Range("N2").Value = 4822.37
$ws.Range("O2").Value = 83
$ws.Range("P2").Value = 3.267716535433071
$ws.Range("Q2").Value = 7066.179999999999
$ws.Range("R2").Value = 540
$ws.Range("S2").Value = 21.25984251968504
$ws.Range("T2").Value = 690
$ws.Range("U2").Value = 27.16535433070866
$ws.Range("V2").Value = 331276.02
$ws.Range("W2").Value = 461
$ws.Range("X2").Value = 18.1496062992126
$ws.Range("Y2").Value = 31
$ws.Range("Z2").Value = 1.220472440944882
$ws.Range("AA2").Value = 1611.75
$ws.Range("AB2").Value = 735
$ws.Range("AC2").Value = 28.93700787401575
$ws.Range("AD2").Value = 2540
$ws.Range("AE2").Value = 2505
$ws.Range("AF2").Value = 98.62204724409449
$ws.Range("AG2").Value = 35
$ws.Range("AH2").Value = 1.377952755905511
$ws.Range("AI2").Value = 125
$ws.Range("AJ2").Value = 175
$ws.Range("AK2").Value = 355
$ws.Range("AL2").Value = 19.08396946564886
$ws.Range("AM2").Value = 26.7175572519084
$ws.Range("AN2").Value = 54.19847328244275
$ws.Range("AO2").Value = 466031.05
$ws.Range("AP2").Value = 87392.91
$ws.Range("AQ2").Value = 29173.04
$ws.Range("AR2").Value = 79.99200991422887
$ws.Range("AS2").Value = 15.00057672799551
$ws.Range("AT2").Value = 5.007413357775614
$ws.Range("AU2").Value = 79.76211175215364
$ws.Range("AV2").Value = 127.9756838905775
$ws.Range("AW2").Value = 200.8143839238498
